$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
$twoItemRows = @(3, 6, 10, 12, 13, 14, 15, 18, 19, 20, 21, 22, 24, 26, 29, 32, 36, 38, 39, 40, 41, 44, 45, 46, 47, 48, 50, 52, 55, 58, 62, 64, 65, 66, 67, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 90, 92, 99, 101, 109, 110, 111, 112, 116, 118, 125, 127, 135, 136, 137, 138, 142, 144, 151, 153)
foreach ($r in $twoItemRows) {
    $ws.Range("G$r").Value = "System, dnasr281@gmail.com"
}

# Rows where "system, backup@backdoor.com, System" -> "backup@backdoor.com, system, System"
$threeItemRows = @(2, 28, 54)
foreach ($r in $threeItemRows) {
    $ws.Range("G$r").Value = "backup@backdoor.com, system, System"
}
